$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 48486.285
$ws.Range("J17").Value = 45329.75
$ws.Range("L17").Value = 135989.25
$ws.Range("N17").Value = -136325.25
$ws.Range("H19").Value = 5952791.5
$ws.Range("I19").Value = 28571618
$ws.Range("J19").Value = 468.8421
$ws.Range("K19").Value = 28571618
$ws.Range("L19").Value = 468.8421
$ws.Range("M19").Value = -28571443
$ws.Range("N19").Value = -818.8421000000001
$ws.Range("H20").Value = 33678
$ws.Range("I20").Value = 4006.6667
$ws.Range("K20").Value = 4006.6667
$ws.Range("M20").Value = -3776.6667
$ws.Range("H34").Value = 5418.6
$ws.Range("I34").Value = 2022
$ws.Range("J34").Value = 7683
$ws.Range("K34").Value = 2022
$ws.Range("L34").Value = 7683
$ws.Range("M34").Value = -1819
$ws.Range("N34").Value = -8089
$ws.Range("H35").Value = 33678
$ws.Range("I35").Value = 4006.6667
$ws.Range("K35").Value = 4006.6667
$ws.Range("M35").Value = -3627.6667
$ws.Range("H36").Value = 5418.6
$ws.Range("I36").Value = 2022
$ws.Range("J36").Value = 7683
$ws.Range("K36").Value = 2022
$ws.Range("L36").Value = 7683
$ws.Range("M36").Value = -1307
$ws.Range("N36").Value = -9113
$ws.Range("H62").Value = 3631.8572
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 4084.6
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 4084.6
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -5332.6
$ws.Range("H65").Value = 3631.8572
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 4084.6
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 20423
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -26663
$ws.Range("H82").Value = 6788.8887
$ws.Range("I82").Value = 5375
$ws.Range("J82").Value = 7920
$ws.Range("K82").Value = 16125
$ws.Range("L82").Value = 23760
$ws.Range("M82").Value = -15719
$ws.Range("N82").Value = -24572
$ws.Range("H85").Value = 6788.8887
$ws.Range("I85").Value = 5375
$ws.Range("J85").Value = 7920
$ws.Range("K85").Value = 16125
$ws.Range("L85").Value = 23760
$ws.Range("M85").Value = -14721
$ws.Range("N85").Value = -26568
$ws.Range("H111").Value = 1393.7778
$ws.Range("I111").Value = 1189.0476
$ws.Range("J111").Value = 2110.3333
$ws.Range("K111").Value = 3567.142800000001
$ws.Range("L111").Value = 6330.999899999999
$ws.Range("M111").Value = -500.1428000000005
$ws.Range("N111").Value = -12464.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 30090.908
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30090.908
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30090.908
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -30766.908
$ws.Range("H79").Value = 30090.908
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30090.908
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30090.908
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -32430.908
$ws.Range("H122").Value = 2375.9556
$ws.Range("I122").Value = 1816.76
$ws.Range("J122").Value = 3074.95
$ws.Range("K122").Value = 5450.28
$ws.Range("L122").Value = 9224.849999999999
$ws.Range("M122").Value = -3000.28
$ws.Range("N122").Value = -14124.85
$ws.Range("H139").Value = 26691.111
$ws.Range("J139").Value = 26691.111
$ws.Range("L139").Value = 26691.111
$ws.Range("N139").Value = -36971.111
$ws.Range("H141").Value = 67785.8
$ws.Range("J141").Value = 67785.8
$ws.Range("L141").Value = 67785.8
$ws.Range("N141").Value = -78145.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 2357.2
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2563.5557
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 2563.5557
$ws.Range("M17").Value = -328
$ws.Range("N17").Value = -2907.5557
$ws.Range("H99").Value = 2651.125
$ws.Range("I99").Value = 1201.5
$ws.Range("K99").Value = 1201.5
$ws.Range("M99").Value = 296.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 24660.3
$ws.Range("I2").Value = 1425
$ws.Range("J2").Value = 40150.5
$ws.Range("K2").Value = 1425
$ws.Range("L2").Value = 40150.5
$ws.Range("M2").Value = -1312
$ws.Range("N2").Value = -40376.5
$ws.Range("H12").Value = 98268.414
$ws.Range("I12").Value = 143315.58
$ws.Range("J12").Value = 35202.4
$ws.Range("K12").Value = 143315.58
$ws.Range("L12").Value = 35202.4
$ws.Range("M12").Value = -143145.58
$ws.Range("N12").Value = -35542.4
$ws.Range("H125").Value = 26666.334
$ws.Range("J125").Value = 26666.334
$ws.Range("L125").Value = 26666.334
$ws.Range("N125").Value = -31586.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1980.6
$ws.Range("I10").Value = 117
$ws.Range("K10").Value = 351
$ws.Range("M10").Value = -212
$ws.Range("H98").Value = 178.66667
$ws.Range("I98").Value = 197.66667
$ws.Range("J98").Value = 159.66667
$ws.Range("K98").Value = 593.00001
$ws.Range("L98").Value = 479.00001
$ws.Range("M98").Value = 904.99999
$ws.Range("N98").Value = -3475.00001
$ws.Range("H99").Value = 2074.75
$ws.Range("I99").Value = 1766.3334
$ws.Range("K99").Value = 5299.0002
$ws.Range("M99").Value = -3053.0002
$ws.Range("H100").Value = 1799
$ws.Range("J100").Value = 2331.6667
$ws.Range("L100").Value = 6995.000100000001
$ws.Range("N100").Value = -8617.000100000001
$ws.Range("H113").Value = 4546310.5
$ws.Range("I113").Value = 100000000
$ws.Range("J113").Value = 896.9524
$ws.Range("K113").Value = 300000000
$ws.Range("L113").Value = 2690.8572
$ws.Range("M113").Value = -299997830
$ws.Range("N113").Value = -7030.8572
$ws.Range("H121").Value = 52150
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 52150
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 156450
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -159070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 68005.39999999999
$ws.Range("J6").Value = 68005.39999999999
$ws.Range("L6").Value = 68005.39999999999
$ws.Range("N6").Value = -68231.39999999999
$ws.Range("H11").Value = 5771765
$ws.Range("J11").Value = 5350000.5
$ws.Range("L11").Value = 5350000.5
$ws.Range("N11").Value = -5350278.5
$ws.Range("H14").Value = 216501
$ws.Range("I14").Value = 1000000
$ws.Range("J14").Value = 20626.25
$ws.Range("K14").Value = 1000000
$ws.Range("L14").Value = 20626.25
$ws.Range("M14").Value = -999832
$ws.Range("N14").Value = -20962.25
$ws.Range("H16").Value = 68005.39999999999
$ws.Range("J16").Value = 68005.39999999999
$ws.Range("L16").Value = 68005.39999999999
$ws.Range("N16").Value = -68505.39999999999
$ws.Range("H31").Value = 3408.8572
$ws.Range("I31").Value = 715.5
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 715.5
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = -423.5
$ws.Range("N31").Value = -7584
$ws.Range("H37").Value = 3408.8572
$ws.Range("I37").Value = 715.5
$ws.Range("J37").Value = 7000
$ws.Range("K37").Value = 715.5
$ws.Range("L37").Value = 7000
$ws.Range("M37").Value = -438.5
$ws.Range("N37").Value = -7554
$ws.Range("H137").Value = 29833.334
$ws.Range("J137").Value = 29833.334
$ws.Range("L137").Value = 29833.334
$ws.Range("N137").Value = -40033.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 380
$ws.Range("I17").Value = 380
$ws.Range("K17").Value = 380
$ws.Range("M17").Value = -210
$ws.Range("H40").Value = 2724.375
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 2970.7144
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 2970.7144
$ws.Range("M40").Value = -864
$ws.Range("N40").Value = -3242.7144
$ws.Range("H106").Value = 32874.75
$ws.Range("J106").Value = 32874.75
$ws.Range("L106").Value = 32874.75
$ws.Range("N106").Value = -35398.75
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 85006
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H81").Value = 892.5
$ws.Range("I81").Value = 713.8461
$ws.Range("K81").Value = 1427.6922
$ws.Range("M81").Value = -366.6922
$ws.Range("H84").Value = 892.5
$ws.Range("I84").Value = 713.8461
$ws.Range("K84").Value = 7138.460999999999
$ws.Range("M84").Value = -1834.460999999999
$ws.Range("H122").Value = 501788.84
$ws.Range("I122").Value = 626392.6
$ws.Range("J122").Value = 3373.75
$ws.Range("K122").Value = 1879177.8
$ws.Range("L122").Value = 10121.25
$ws.Range("M122").Value = -1876727.8
$ws.Range("N122").Value = -15021.25
$ws.Range("H132").Value = 156958.27
$ws.Range("I132").Value = 197203.7
$ws.Range("J132").Value = 10349.857
$ws.Range("K132").Value = 591611.1000000001
$ws.Range("L132").Value = 31049.571
$ws.Range("M132").Value = -589081.1000000001
$ws.Range("N132").Value = -36109.571
